$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select and delete row 1, shifting the existing data (rows 2-3) up by one row
# so it ends up occupying rows 1-2, matching the target layout.
$ws.Rows.Item(1).Select()
$ws.Rows.Item(1).Delete()
